$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are stored as text so formatted strings
# like '28.80' or '0.0901' keep their exact digits instead of being
# reinterpreted as numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.887.05'
$ws.Range("E2").Value = '  +0.03%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.635.42'
$ws.Range("E3").Value = '  +0.81%  '
$ws.Range("E4").Value = '  +0.76%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '215.41'
$ws.Range("E5").Value = '  +0.89%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.520'
$ws.Range("E6").Value = '  +0.11%  '
$ws.Range("E7").Value = '  +0.84%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '28.80'
$ws.Range("E8").Value = '  -1.15%  '
$ws.Range("E9").Value = '  +0.70%  '
$ws.Range("E10").Value = '  +0.39%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0901'
$ws.Range("E11").Value = '  -1.13%  '
$ws.Range("B12").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C12").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.869.90'
$ws.Range("E12").Value = '  +0.89%  '
$ws.Range("B13").Value = 'WrappedEther'
$ws.Range("C13").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.630.81'
$ws.Range("E13").Value = '  +0.61%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.590'
$ws.Range("E14").Value = '  +3.86%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '9.43'
$ws.Range("E15").Value = '  +6.63%  '
$ws.Range("B16").Value = 'WrappedBTC'
$ws.Range("C16").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '29.910.46'
$ws.Range("E16").Value = '  +0.06%  '
$ws.Range("B17").Value = 'Polkadot'
$ws.Range("C17").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.85'
$ws.Range("E17").Value = '  -1.47%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '64.62'
$ws.Range("E18").Value = '  +0.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '240.64'
$ws.Range("E19").Value = '  -0.58%  '
$ws.Range("E21").Value = '  +0.65%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.89'
$ws.Range("E22").Value = '  +3.11%  '
$ws.Range("E23").Value = '  +0.98%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.17'
$ws.Range("E24").Value = '  +2.81%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '157.24'
$ws.Range("E25").Value = '  +0.95%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '15.53'
$ws.Range("E26").Value = '  -0.62%  '
$ws.Range("E27").Value = '  -0.67%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.63'
$ws.Range("E28").Value = '  +0.66%  '
$ws.Range("E29").Value = '  +0.65%  '
$ws.Range("E30").Value = '  +0.74%  '
$ws.Range("E31").Value = '  -0.79%  '
$ws.Range("E32").Value = '  +1.08%  '
$ws.Range("E33").Value = '  -0.68%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.425.93'
$ws.Range("E34").Value = '  +0.59%  '
$ws.Range("E35").Value = '  +2.85%  '
$ws.Range("E36").Value = '  -0.95%  '
$ws.Range("E37").Value = '  -3.21%  '
$ws.Range("E38").Value = '  +1.58%  '
$ws.Range("E39").Value = '  +0.22%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '76.18'
$ws.Range("E40").Value = '  +10.04%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.561'
$ws.Range("E41").Value = '  +0.98%  '
$ws.Range("E42").Value = '  +0.56%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0500'
$ws.Range("E43").Value = '  -0.55%  '
$ws.Range("E44").Value = '  -0.47%  '
$ws.Range("E45").Value = '  +0.76%  '
$ws.Range("E46").Value = '  -1.42%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '5.38'
$ws.Range("E47").Value = '  -0.82%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.777.27'
$ws.Range("E48").Value = '  +0.85%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '48.87'
$ws.Range("E49").Value = '  -8.49%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '92.97'
$ws.Range("E50").Value = '  +5.05%  '
$ws.Range("E51").Value = '  +0.42%  '
